# custom accuracy + 데이터 1000개
# Round the last data row (row 5) to a custom accuracy (2 decimal places)
# and drop the trailing row (row 6) that no longer belongs to the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "custom accuracy" (2 decimal places) rounding to row 5, B:AH
$ws.Range("B5").Value  = 15.85
$ws.Range("C5").Value  = 11.57
$ws.Range("D5").Value  = 1.09
$ws.Range("E5").Value  = 34.4
$ws.Range("F5").Value  = 28.11
$ws.Range("G5").Value  = 12.48
$ws.Range("H5").Value  = 48.14
$ws.Range("I5").Value  = 19.2
$ws.Range("J5").Value  = 8.46
$ws.Range("K5").Value  = 12.52
$ws.Range("L5").Value  = 13.82
$ws.Range("M5").Value  = 14.52
$ws.Range("N5").Value  = 3.99
$ws.Range("O5").Value  = 12.41
$ws.Range("P5").Value  = 17.6
$ws.Range("Q5").Value  = 10.54
$ws.Range("R5").Value  = 0.81
$ws.Range("S5").Value  = 0.7
$ws.Range("T5").Value  = 181.52
$ws.Range("U5").Value  = 34.7
$ws.Range("V5").Value  = 11.45
$ws.Range("W5").Value  = 23.24
$ws.Range("X5").Value  = 12.18
$ws.Range("Y5").Value  = 1.92
$ws.Range("Z5").Value  = 23.38
$ws.Range("AA5").Value = 10.12
$ws.Range("AB5").Value = 9.03
$ws.Range("AC5").Value = 10.6
$ws.Range("AD5").Value = 14.43
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 43.68
$ws.Range("AG5").Value = 6.42
$ws.Range("AH5").Value = 14.32

# Remove the now-obsolete last row (row 6) of the sample data
$ws.Rows("6").Delete()
